# Append new daily NAV/basket rows (2024-08-28 .. 2024-09-23) to Sheet1.
# Columns: A=Date (text), C..G=stock prices, H=Basket Value, I=daily return, J=NAV
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(594, "2024-08-28", 1074.550048828125, 700.0999755859375, 260.4500122070312, 491.7000122070312, 1535.849975585938, 29293.15045166016, -0.009101137231307158, 274.1078956819022),
    @(595, "2024-08-29", 1121.650024414062, 689.1500244140625, 253.8000030517578, 490.2000122070312, 1520.25, 29353.35063171387, 0.002055094079179154, 274.6712131953744),
    @(596, "2024-08-30", 1111.349975585938, 683.3499755859375, 256.4500122070312, 492.8999938964844, 1521, 29291.19976806641, -0.002117334556700047, 274.0896423439451),
    @(597, "2024-09-02", 1092.650024414062, 715.0499877929688, 251.3500061035156, 490.5, 1505.25, 29322.30020141602, 0.00106176713811209, 274.3806617190828),
    @(598, "2024-09-03", 1085.099975585938, 710.7999877929688, 251, 488.8500061035156, 1509, 29201.39971923828, -0.004123158188384413, 273.2493468469814),
    @(599, "2024-09-04", 1080.449951171875, 722.4000244140625, 250.5, 484.1499938964844, 1488.099975585938, 29200.89978027344, -0.00001712037675078922, 273.2446687152164),
    @(600, "2024-09-05", 1069.150024414062, 733.8499755859375, 251.1499938964844, 495.6499938964844, 1447.599975585938, 29262.3996887207, 0.002106096350113556, 273.8201483146856),
    @(601, "2024-09-06", 1049.349975585938, 718.9000244140625, 247.8000030517578, 483, 1418.050048828125, 28702.20024108887, -0.01914400232349252, 268.5781347591301),
    @(602, "2024-09-09", 1038.699951171875, 700.1500244140625, 243.8999938964844, 474.75, 1411.849975585938, 28242.6496887207, -0.01601098691069303, 264.2779337590034),
    @(603, "2024-09-10", 1035.800048828125, 713.4000244140625, 248.25, 478.7999877929688, 1424.449951171875, 28522.85046386719, 0.009921192885042528, 266.8998861150869),
    @(604, "2024-09-11", 976.2999877929688, 725.4000244140625, 241.5500030517578, 472.2000122070312, 1399.599975585938, 27922.25028991699, -0.02105680758348599, 261.2798265691072),
    @(605, "2024-09-12", 986.1500244140625, 726.0499877929688, 246.1499938964844, 479.8500061035156, 1403.150024414062, 28182.30001831055, 0.009313351384414074, 263.7132174036041),
    @(606, "2024-09-13", 992.0999755859375, 724.25, 245.6499938964844, 485.3999938964844, 1410.949951171875, 28258.59951782227, 0.002707355306775728, 264.4271827822087),
    @(607, "2024-09-16", 988.4000244140625, 733.6500244140625, 243.8000030517578, 489.9500122070312, 1404.550048828125, 28313.45072937012, 0.00194104493795801, 264.9404478268066),
    @(608, "2024-09-17", 974.9500122070312, 745.4000244140625, 240.8000030517578, 482.2999877929688, 1400.25, 28196.30033874512, -0.004137623200533361, 263.8442240831187),
    @(609, "2024-09-18", 962.0499877929688, 717.5499877929688, 235.9499969482422, 471.75, 1391.300048828125, 27572.89979553223, -0.02210930284198537, 258.0108122297563),
    @(610, "2024-09-19", 967, 728.5, 237.5500030517578, 459.9500122070312, 1374.150024414062, 27641.0502166748, 0.002471645044516532, 258.6485233752357),
    @(611, "2024-09-20", 970.8499755859375, 748.3499755859375, 237.8500061035156, 466.2999877929688, 1380.550048828125, 27960.69967651367, 0.01156430227264067, 261.6396130819191),
    @(612, "2024-09-23", 971.7999877929688, 750.2000122070312, 236.4499969482422, 471.1499938964844, 1375.400024414062, 27984.94996643066, 0.0008672991090191444, 261.8665328852292)
)

foreach ($r in $newRows) {
    $rowNum = $r[0]

    # Column A holds the date as plain text (matches existing rows), so force
    # text formatting before assigning to avoid Excel's automatic date parsing.
    $aCell = $ws.Cells.Item($rowNum, 1)
    $aCell.NumberFormat = "@"
    $aCell.Value = $r[1]

    $ws.Cells.Item($rowNum, 3).Value  = $r[2]   # C
    $ws.Cells.Item($rowNum, 4).Value  = $r[3]   # D
    $ws.Cells.Item($rowNum, 5).Value  = $r[4]   # E
    $ws.Cells.Item($rowNum, 6).Value  = $r[5]   # F
    $ws.Cells.Item($rowNum, 7).Value  = $r[6]   # G
    $ws.Cells.Item($rowNum, 8).Value  = $r[7]   # H
    $ws.Cells.Item($rowNum, 9).Value  = $r[8]   # I
    $ws.Cells.Item($rowNum, 10).Value = $r[9]   # J
}

Write-Host "Added $($newRows.Count) rows (594-612) to Sheet1."
